# USU-08 ready and Refactor by providers backend
#
# Marks the "USU-02" (Registro con Google/Facebook) and "USU-04"
# (Inicio de Sesion con Google/Facebook) user-story rows in the first
# table as ready by colouring all of their text green (RGB 146,208,80
# / hex 92D050), matching the styling already used on the other
# completed rows (USU-01, USU-03, etc.) in that table.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# RGB(146, 208, 80) == hex 92D050, encoded as a Word OLE color
# (0x00BBGGRR) the way VBA's RGB() helper would produce it.
$readyGreen = 146 + (208 * 256) + (80 * 65536)

# Row 3 = USU-02 ("Registro con Google/Facebook"), Row 5 = USU-04
# ("Inicio de Sesion con Google/Facebook"). Colour every cell (Codigo,
# Historia de Usuario, Criterios de Aceptacion) in both rows.
$readyRows = @(3, 5)

foreach ($rowIndex in $readyRows) {
    $row = $table.Rows.Item($rowIndex)
    for ($colIndex = 1; $colIndex -le $row.Cells.Count; $colIndex++) {
        $cell = $table.Cell($rowIndex, $colIndex)
        $cell.Range.Font.Color = $readyGreen
    }
}
